# Weekly refresh of fruit/vegetable prices: the rows 2..49 of the sheet keep
# the same set of (Fecha, Volumen, Precio minimo, Precio maximo,
# Precio promedio ponderado, Precio $/Kg) records, but they get reshuffled
# onto different rows. Row 22 is unchanged.
#
# Mapping below means: newRow[r] = oldRow[mapping[r]]  (values only, for
# columns D, J, K, L, M, P). All other columns stay put.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$mapping = @{
    2  = 34
    3  = 29
    4  = 7
    5  = 25
    6  = 40
    7  = 18
    8  = 35
    9  = 26
    10 = 2
    11 = 3
    12 = 14
    13 = 21
    14 = 32
    15 = 24
    16 = 11
    17 = 27
    18 = 49
    19 = 31
    20 = 42
    21 = 30
    22 = 22
    23 = 17
    24 = 4
    25 = 9
    26 = 10
    27 = 33
    28 = 44
    29 = 28
    30 = 12
    31 = 41
    32 = 15
    33 = 48
    34 = 46
    35 = 5
    36 = 45
    37 = 8
    38 = 36
    39 = 20
    40 = 19
    41 = 38
    42 = 13
    43 = 16
    44 = 6
    45 = 23
    46 = 47
    47 = 43
    48 = 39
    49 = 37
}

$cols = @("D", "J", "K", "L", "M", "P")

# Snapshot the original values for the columns that move, indexed by row.
$snapshot = @{}
foreach ($row in $mapping.Keys) {
    $rowValues = @{}
    foreach ($col in $cols) {
        $rowValues[$col] = $ws.Range("$col$row").Value2
    }
    $snapshot[$row] = $rowValues
}

# Write back according to the mapping, using the snapshot as the source so
# that earlier writes never disturb later reads.
foreach ($row in $mapping.Keys) {
    $srcRow = $mapping[$row]
    $srcValues = $snapshot[$srcRow]
    foreach ($col in $cols) {
        $ws.Range("$col$row").Value = $srcValues[$col]
    }
}
